$d = $word.ActiveDocument

# Delete the three paragraphs: "Create model of cafes", "Add café to DB",
# "Create superuser" - they appear right after "TODO" and before
# "Create form for adding new places".

foreach ($target in @("Create model of cafes", "Add café to DB", "Create superuser")) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd("`r") -eq $target) {
            $p.Range.Delete()
            break
        }
    }
}
